# Updated symbol list on Mon Dec 19 07:28:25 UTC 2022 with GitHub Actions
#
# The "Price" column (D) and part of the "Volume(1h)" column (E) hold
# text values (e.g. "247.80"), not real numbers, so each numeric-looking
# price is written with a leading apostrophe to force a text entry and
# preserve the exact formatting (trailing zeros etc.) instead of letting
# Excel coerce it into a floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'247.62"
$ws.Range("D3").Value  = "'21.76"
$ws.Range("D4").Value  = "'5.454"
$ws.Range("D5").Value  = "'0.05679"
$ws.Range("D7").Value  = "'0.8014"
$ws.Range("D8").Value  = "'1.037"
$ws.Range("D9").Value  = "'0.1451"
$ws.Range("D10").Value = "'0.07247"
$ws.Range("D12").Value = "'0.02936"
$ws.Range("D13").Value = "'0.09284"
$ws.Range("D14").Value = "'0.001659"
$ws.Range("D15").Value = "'3.206"
$ws.Range("E15").Value = "14MCDexMCBBestin24h"
$ws.Range("D16").Value = "'0.04719"
$ws.Range("E17").Value = "16OneONEWorstin24h"
$ws.Range("D18").Value = "'0.006349"
$ws.Range("D22").Value = "'0.0003200"
$ws.Range("D23").Value = "'3.804"
$ws.Range("D24").Value = "'6.423"
$ws.Range("D25").Value = "'2.125"
$ws.Range("D27").Value = "'0.1298"
$ws.Range("D40").Value = "'0.04087"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D44").Value = "'0.008948"
$ws.Range("D45").Value = "'0.00005847"
$ws.Range("D47").Value = "'0.7853"
$ws.Range("D48").Value = "'0.01046"
$ws.Range("E48").Value = "47BOLOBOLO"
